$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcc77e613f3649b452caacb6a8f2a54d142985ae/e2e/b6ffff3b-a984-4897-944a-f55dee0f69ba.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c055d25e42a4c8b905b970af3fb983691772e7f0/e2e/b6ffff3b-a984-4897-944a-f55dee0f69ba.md."

# --- Overview sheet: row for b6ffff3b-a984-4897-944a-f55dee0f69ba.md is now fully handed off ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-24 20:50:21"

# --- zh-cn sheet: same row (row 3) reflects the new handoff ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-24 20:50:00"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet: same row (row 3) reflects the new handoff ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-24 20:50:21"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
